# The edit re-orders the observation records that occupy rows 2-10 of the
# active sheet (each row is uniquely identified by its "Id" in column A).
# Every column of a given record moves together - this script reproduces
# that end state by writing the new value combination directly into each
# target row/column, and clearing the few cells (column AC, AJ, AK, AO)
# that become empty for a given row after the reshuffle.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (was Id 80448771, now Id 80448769)
$ws.Range("A2").Value = 80448769
$ws.Range("B2").Value = 77506
$ws.Range("E2").Value = 6425
$ws.Range("F2").Value = "Garnlav"
$ws.Range("G2").Value = "Alectoria sarmentosa"
$ws.Range("H2").Value = "(Ach.) Ach."
$ws.Range("Q2").Value = 422991.0759451608
$ws.Range("R2").Value = 6752021.173145968
$ws.Range("AC2").Value = "Rikligt, hkb"

# Row 3 (was Id 80448772, now Id 80448775)
$ws.Range("A3").Value = 80448775
$ws.Range("B3").Value = 77506
$ws.Range("E3").Value = 6425
$ws.Range("F3").Value = "Garnlav"
$ws.Range("G3").Value = "Alectoria sarmentosa"
$ws.Range("H3").Value = "(Ach.) Ach."
$ws.Range("Q3").Value = 423036.1594514723
$ws.Range("R3").Value = 6752009.000504656
$ws.Range("AC3").Value = "Rikligt, hkb"

# Row 4 (was Id 80448769, now Id 80448777)
$ws.Range("A4").Value = 80448777
$ws.Range("Q4").Value = 423115.1561234437
$ws.Range("R4").Value = 6752009.239606674
$ws.Range("AC4").ClearContents()

# Row 5 (was Id 80448775, now Id 80448771)
$ws.Range("A5").Value = 80448771
$ws.Range("B5").Value = 73693
$ws.Range("E5").Value = 6440
$ws.Range("F5").Value = "Vitgrynig nållav"
$ws.Range("G5").Value = "Chaenotheca subroscida"
$ws.Range("H5").Value = "(Eitner) Zahlbr."
$ws.Range("Q5").Value = 423289.9356373397
$ws.Range("R5").Value = 6752041.978126496
$ws.Range("AC5").ClearContents()

# Row 6 (was Id 80448777, now Id 80448779)
$ws.Range("A6").Value = 80448779
$ws.Range("Q6").Value = 422962.8083476268
$ws.Range("R6").Value = 6752021.785183201
$ws.Range("AC6").Value = "Rikligt"

# Row 7 (was Id 80448779, now Id 80448780)
$ws.Range("A7").Value = 80448780
$ws.Range("Q7").Value = 423056.1482692101
$ws.Range("R7").Value = 6751963.779848268
$ws.Range("AC7").Value = "Spritt"

# Row 8 (was Id 80448780, now Id 80448772)
$ws.Range("A8").Value = 80448772
$ws.Range("B8").Value = 81236
$ws.Range("E8").Value = 1312
$ws.Range("F8").Value = "Gammelgransskål"
$ws.Range("G8").Value = "Pseudographis pinicola"
$ws.Range("H8").Value = "(Nyl.) Rehm"
$ws.Range("Q8").Value = 423289.9356373397
$ws.Range("R8").Value = 6752041.978126496
$ws.Range("AC8").ClearContents()

# Row 9 (was Id 80448778, now Id 80448773)
$ws.Range("A9").Value = 80448773
$ws.Range("B9").Value = 77506
$ws.Range("C9").Value = "Ovaliderad"
$ws.Range("E9").Value = 6425
$ws.Range("F9").Value = "Garnlav"
$ws.Range("G9").Value = "Alectoria sarmentosa"
$ws.Range("H9").Value = "(Ach.) Ach."
$ws.Range("Q9").Value = 422635.9957601223
$ws.Range("R9").Value = 6751949.037152009
$ws.Range("AJ9").Value = "vanlig tall"
$ws.Range("AK9").Value = "Pinus sylvestris var. sylvestris"
$ws.Range("AO9").Value = "Pinus sylvestris var. sylvestris"

# Row 10 (was Id 80448773, now Id 80448778)
$ws.Range("A10").Value = 80448778
$ws.Range("B10").Value = 56395
$ws.Range("C10").Value = "Godkänd baserat på observatörens uppgifter"
$ws.Range("E10").Value = 100109
$ws.Range("F10").Value = "Tretåig hackspett"
$ws.Range("G10").Value = "Picoides tridactylus"
$ws.Range("H10").Value = "(Linnaeus, 1758)"
$ws.Range("Q10").Value = 423115.1561234437
$ws.Range("R10").Value = 6752009.239606674
$ws.Range("AJ10").ClearContents()
$ws.Range("AK10").ClearContents()
$ws.Range("AO10").ClearContents()
